$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 5360, 45890),
    @(3, 5300, 45890.01041666666),
    @(4, 5260, 45890.02083333334),
    @(5, 5220, 45890.03125),
    @(6, 5180, 45890.04166666666),
    @(7, 5150, 45890.05208333334),
    @(8, 5140, 45890.0625),
    @(9, 5130, 45890.07291666666),
    @(10, 5110, 45890.08333333334),
    @(11, 5100, 45890.09375),
    @(12, 5100, 45890.10416666666),
    @(13, 5100, 45890.11458333334),
    @(14, 5100, 45890.125),
    @(15, 5100, 45890.13541666666),
    @(16, 5100, 45890.14583333334),
    @(17, 5110, 45890.15625),
    @(18, 5150, 45890.16666666666),
    @(19, 5190, 45890.17708333334),
    @(20, 5240, 45890.1875),
    @(21, 5280, 45890.19791666666),
    @(22, 5340, 45890.20833333334),
    @(23, 5400, 45890.21875),
    @(24, 5490, 45890.22916666666),
    @(25, 5580, 45890.23958333334),
    @(26, 5750, 45890.25),
    @(27, 5870, 45890.26041666666),
    @(28, 5930, 45890.27083333334),
    @(29, 5960, 45890.28125),
    @(30, 5980, 45890.29166666666),
    @(31, 5980, 45890.30208333334),
    @(32, 5970, 45890.3125),
    @(33, 5930, 45890.32291666666),
    @(34, 5820, 45890.33333333334),
    @(35, 5730, 45890.34375),
    @(36, 5660, 45890.35416666666),
    @(37, 5580, 45890.36458333334),
    @(38, 5500, 45890.375),
    @(39, 5420, 45890.38541666666),
    @(40, 5350, 45890.39583333334),
    @(41, 5300, 45890.40625),
    @(42, 5230, 45890.41666666666),
    @(43, 5200, 45890.42708333334),
    @(44, 5190, 45890.4375),
    @(45, 5180, 45890.44791666666),
    @(46, 5170, 45890.45833333334),
    @(47, 5170, 45890.46875),
    @(48, 5170, 45890.47916666666),
    @(49, 5170, 45890.48958333334),
    @(50, 5190, 45890.5),
    @(51, 5210, 45890.51041666666),
    @(52, 5230, 45890.52083333334),
    @(53, 5260, 45890.53125),
    @(54, 5310, 45890.54166666666),
    @(55, 5360, 45890.55208333334),
    @(56, 5400, 45890.5625),
    @(57, 5450, 45890.57291666666),
    @(58, 5510, 45890.58333333334),
    @(59, 5570, 45890.59375),
    @(60, 5630, 45890.60416666666),
    @(61, 5690, 45890.61458333334),
    @(62, 5780, 45890.625),
    @(63, 5860, 45890.63541666666),
    @(64, 5960, 45890.64583333334),
    @(65, 6060, 45890.65625),
    @(66, 6160, 45890.66666666666),
    @(67, 6250, 45890.67708333334),
    @(68, 6350, 45890.6875),
    @(69, 6440, 45890.69791666666),
    @(70, 6570, 45890.70833333334),
    @(71, 6670, 45890.71875),
    @(72, 6780, 45890.72916666666),
    @(73, 6880, 45890.73958333334),
    @(74, 6960, 45890.75),
    @(75, 7060, 45890.76041666666),
    @(76, 7120, 45890.77083333334),
    @(77, 7170, 45890.78125),
    @(78, 7200, 45890.79166666666),
    @(79, 7230, 45890.80208333334),
    @(80, 7260, 45890.8125),
    @(81, 7310, 45890.82291666666),
    @(82, 7320, 45890.83333333334),
    @(83, 7290, 45890.84375),
    @(84, 7240, 45890.85416666666),
    @(85, 7130, 45890.86458333334),
    @(86, 6940, 45890.875),
    @(87, 6780, 45890.88541666666),
    @(88, 6630, 45890.89583333334),
    @(89, 6470, 45890.90625),
    @(90, 6290, 45890.91666666666),
    @(91, 6140, 45890.92708333334),
    @(92, 6030, 45890.9375),
    @(93, 5920, 45890.94791666666),
    @(94, 5740, 45890.95833333334),
    @(95, 5680, 45890.96875),
    @(96, 5630, 45890.97916666666),
    @(97, 5570, 45890.98958333334)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}

Write-Host "Done updating rows"